# Weekly fruit/vegetable price update: insert a new daily reading as row 18
# (Feria Lagunitas de Puerto Montt - Haba), shifting the existing rows
# 18-42 down to 19-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18, pushing everything below it down
# by one (this also grows the used range from R42 to R43, matching the
# Excel-authored diff).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(18, 1).Value = 4
$ws.Cells.Item(18, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(18, 3).Value = "Los Lagos"
$ws.Cells.Item(18, 4).Value = 44483
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(18, 6).Value = 100112026
$ws.Cells.Item(18, 7).Value = "Haba"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 70
$ws.Cells.Item(18, 11).Value = 11000
$ws.Cells.Item(18, 12).Value = 11000
$ws.Cells.Item(18, 13).Value = 11000
$ws.Cells.Item(18, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 440
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
